$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Update the date in A1 (was 45310 / 2024-01-19, now 45311 / 2024-01-20)
$ws.Range("A1").Value = 45311

# Update the prices in column D (fix bug exceeded request in google drive)
$ws.Range("D34").Value = 110
$ws.Range("D35").Value = 117
$ws.Range("D36").Value = 134.6
$ws.Range("D37").Value = 114.7
$ws.Range("D38").Value = 128.7
$ws.Range("D39").Value = 160.3
